# "Added Tab with resources and URLs"
#
# 1. Rename the original sheet to "List of Projects".
# 2. Add a new "Sources" sheet right after it, listing project-idea
#    resources (URL + description), and make it the active tab.
# 3. Sheet1's old selection (D5) is replaced by a later one (C39) and it
#    is no longer the selected/active tab.

$wb = $excel.ActiveWorkbook

# --- Rename the existing sheet -------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "List of Projects"

# --- Insert the new "Sources" sheet right after it -----------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sources"

# --- Header row -----------------------------------------------------------------
$ws2.Range("C4").Value = "URL"
$ws2.Range("D4").Value = "Description"

$header = $ws2.Range("C4:D4")
$header.VerticalAlignment = -4108   # xlCenter
$header.HorizontalAlignment = -4108 # xlCenter
$header.Font.Bold = $true

# --- Data rows (order mirrors how the author originally typed them) ------------
$ws2.Range("C5").Value = "https://github.com/JSmolins/Martyrs-Mega-List"
$ws2.Range("C6").Value = "https://github.com/karan/Projects#classic-algorithms"
$ws2.Range("D6").Value = "Karan Mega Project List"
$ws2.Range("D5").Value = "Martyr2's Mega Project List"
$ws2.Range("C7").Value = "https://www.youtube.com/watch?v=Bj6N0pEVC-I"
$ws2.Range("D7").Value = "Chess Stepping Off Point"
$ws2.Range("C8").Value = "https://www.reddit.com/r/dailyprogrammer/"
$ws2.Range("D8").Value = "Subreddit for Daily Challenges"

# --- Data row alignment (left / vertical-center) --------------------------------
$data = $ws2.Range("C5:D8")
$data.VerticalAlignment = -4108    # xlCenter
$data.HorizontalAlignment = -4131  # xlLeft

# --- Column widths ----------------------------------------------------------------
$ws2.Columns("C").ColumnWidth = 50.42
$ws2.Columns("D").ColumnWidth = 24.92

# --- Selections / active tab -----------------------------------------------------
$ws1.Range("C39").Select() | Out-Null
$ws2.Range("C16").Select() | Out-Null
$ws2.Activate() | Out-Null
